# Update the Wnt2-Fzd5 LR-pair sheet with refreshed TPM-based NATMI output.
# The previous data only covered the "FAPs" sending cluster against the
# 3 target clusters (ECs, FAPs, MuSCs). The refreshed script output adds a
# second sending cluster ("ECs") with its own 3 target rows, and replaces
# every numeric metric with newly computed TPM values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A..T =
# Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# Ligand-expressing cells, Ligand detection rate,
# Ligand average expression value, Ligand total expression value,
# Ligand derived specificity of average expression value,
# Ligand derived specificity of total expression value,
# Receptor-expressing cells, Receptor detection rate,
# Receptor average expression value, Receptor total expression value,
# Receptor derived specificity of average expression value,
# Receptor derived specificity of total expression value,
# Edge average expression weight, Edge total expression weight,
# Edge average expression derived specificity, Edge total expression derived specificity
$data = @(
  @("ECs",  "Wnt2", "Fzd5", "ECs",   1, 0.3333333333333333, 0.01070233333333333, 0.032107,
    0.004227647500550067, 0.004227647500550067, 3, 1,
    2.133443333333334, 6.40033, 0.2605947899689859, 0.2605947899689859,
    0.02283282170111111, 0.20549539531, 0.001101702912468753, 0.001101702912468753),
  @("ECs",  "Wnt2", "Fzd5", "FAPs",  1, 0.3333333333333333, 0.01070233333333333, 0.032107,
    0.004227647500550067, 0.004227647500550067, 3, 1,
    4.264793333333333, 12.79438, 0.5209338844846115, 0.5209338844846116,
    0.04564323985111111, 0.41078915866, 0.002202324834693205, 0.002202324834693205),
  @("ECs",  "Wnt2", "Fzd5", "MuSCs", 1, 0.3333333333333333, 0.01070233333333333, 0.032107,
    0.004227647500550067, 0.004227647500550067, 3, 1,
    1.788586, 5.365758, 0.2184713255464024, 0.2184713255464024,
    0.01914204356733333, 0.172278392106, 0.0009236197533881081, 0.0009236197533881083),
  @("FAPs", "Wnt2", "Fzd5", "ECs",   3, 1, 2.520808, 7.562424,
    0.99577235249945, 0.99577235249945, 3, 1,
    2.133443333333334, 6.40033, 0.2605947899689859, 0.2605947899689859,
    5.378001022213335, 48.40200919992, 0.2594930870565172, 0.2594930870565172),
  @("FAPs", "Wnt2", "Fzd5", "FAPs",  3, 1, 2.520808, 7.562424,
    0.99577235249945, 0.99577235249945, 3, 1,
    4.264793333333333, 12.79438, 0.5209338844846115, 0.5209338844846116,
    10.75072515301333, 96.75652637712, 0.5187315596499184, 0.5187315596499185),
  @("FAPs", "Wnt2", "Fzd5", "MuSCs", 3, 1, 2.520808, 7.562424,
    0.99577235249945, 0.99577235249945, 3, 1,
    1.788586, 5.365758, 0.2184713255464024, 0.2184713255464024,
    4.508681897488, 40.578137077392, 0.2175477057930143, 0.2175477057930143)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}
